# Applies the "Trade #16 closed" update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet - update headline metrics
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.28   # Current Capital
$summary.Range("B4").Value = -0.72     # Total P&L $
$summary.Range("B5").Value = -0.9      # Total P&L %
$summary.Range("B6").Value = 16        # Total Trades
$summary.Range("B7").Value = 5         # Winning Trades
$summary.Range("B9").Value = 31.25     # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - update the MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.28
$status.Range("D4").Value = 16
$status.Range("E4").Value = -0.72
$status.Range("F4").Value = -0.72
$status.Range("G4").Value = 31.25

# ---------------------------------------------------------------------------
# 3. Append the newly closed trade (#16) to the "All Trades" and
#    "MarketMaking" trade log sheets (both share the same layout).
#    Note: the Date column ("2026-02-17") is prefixed with a leading
#    apostrophe so Excel stores it as literal text instead of auto-
#    converting it into a date serial number (matching the other rows,
#    which all hold the date as plain text).
# ---------------------------------------------------------------------------
$newTradeRow = @(16, "'2026-02-17", "13:17:33", "MarketMaking", "UP", 0.2, 0.23, "CLOSED", 15, 0.03, 99.28, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 17
    $col = 1
    foreach ($val in $newTradeRow) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
}
